# Thong_tin_hang_trong_kho.xlsx - bo sung cot "So luong yeu cau" (orderAmount)
#
# Layout change: a new column is inserted right after the existing
# "So luong" (amountValue) column. The old "So luong" column is
# relabelled "So luong ton kho" and keeps ${item.amountValue}; the new
# column becomes "So luong yeu cau" / ${item.orderAmountValue}; the
# columns that used to follow (issueAmountValue / goodsUnitName /
# changeDate) simply shift one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural edit: insert a new column G (old F = "So luong" stays,
#     old G/H/I shift right to H/I/J) -----------------------------------
$ws.Columns.Item(7).Insert()

# --- header row (row 3) text updates -----------------------------------
$ws.Range("G3").Value = 'Số lượng yêu cầu'
$ws.Range("F3").Value = 'Số lượng tồn kho'

# --- template placeholder row (row 5) ------------------------------------
$ws.Range("G5").Value = '${item.orderAmountValue}'

# --- column widths (characters) -----------------------------------------
$ws.Columns.Item(1).ColumnWidth = 4.709635416666667
$ws.Columns.Item(2).ColumnWidth = 12.709635416666666
$ws.Columns.Item(3).ColumnWidth = 34.166666666666664
$ws.Columns.Item(4).ColumnWidth = 11.709635416666666
$ws.Columns.Item(5).ColumnWidth = 31.256510416666668
$ws.Columns.Item(6).ColumnWidth = 11.256510416666666
$ws.Columns.Item(7).ColumnWidth = 11.346354166666666
$ws.Columns.Item(8).ColumnWidth = 12.346354166666666
$ws.Columns.Item(9).ColumnWidth = 8.346354166666666
$ws.Columns.Item(10).ColumnWidth = 15.619791666666666

# --- row heights ----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 23.5
$ws.Rows.Item(2).RowHeight = 15.65
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 24.65
$ws.Rows.Item(6).RowHeight = 15

# --- selection / active cell ----------------------------------------------
$ws.Range("E14").Select()

Write-Host "orderAmount column inserted"
